$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 182
$ws.Range("I2").Value = 488
$ws.Range("J2").Value = 2027
$ws.Range("K2").Value = 13
$ws.Range("L2").Value = 557
$ws.Range("M2").Value = 38
$ws.Range("N2").Value = 349
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 30
$ws.Range("S2").Value = 202
$ws.Range("T2").Value = 318
$ws.Range("V2").Value = 3249
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 3042
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 45
$ws.Range("AA2").Value = 17
